# Applies the "Updated cryptos list" crypto price/volume refresh.
# Column D ("Price") cells that look like plain numbers (e.g. "1.00", "571.14")
# are written with a leading apostrophe so Excel stores them as literal text
# (matching the original inlineStr cells) instead of silently coercing them
# to numeric values. Two D39/D40 rows also swap their Coin/Link/Price/Volume
# content (InjectiveProtocol <-> PEPE changed rank order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.860.26'
$ws.Range("E2").Value = '  -3.20%  '
$ws.Range("D3").Value = '3.435.77'
$ws.Range("E3").Value = '  -2.91%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'" + '571.14'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = "'" + '174.64'
$ws.Range("E6").Value = '  -8.20%  '
$ws.Range("D7").Value = "'" + '0.622'
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("E10").Value = '  +4.85%  '
$ws.Range("D11").Value = "'" + '54.84'
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = "'" + '0.0000272'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = "'" + '9.11'
$ws.Range("E13").Value = '  -3.42%  '
$ws.Range("D14").Value = '3.983.21'
$ws.Range("E14").Value = '  -3.02%  '
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").Value = '3.441.88'
$ws.Range("E16").Value = '  -2.87%  '
$ws.Range("D17").Value = "'" + '18.07'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '64.856.11'
$ws.Range("E18").Value = '  -3.28%  '
$ws.Range("D19").Value = "'" + '11.83'
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").Value = "'" + '0.987'
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("D21").Value = "'" + '407.21'
$ws.Range("E21").Value = '  -6.00%  '
$ws.Range("D22").Value = "'" + '4.18'
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").Value = "'" + '4.35'
$ws.Range("E23").Value = '  +5.01%  '
$ws.Range("D24").Value = "'" + '83.51'
$ws.Range("E24").Value = '  -2.15%  '
$ws.Range("D25").Value = "'" + '13.21'
$ws.Range("E25").Value = '  +7.76%  '
$ws.Range("D26").Value = "'" + '10.80'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("E27").Value = '  -3.63%  '
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").Value = "'" + '8.94'
$ws.Range("E29").Value = '  -2.48%  '
$ws.Range("D30").Value = "'" + '29.82'
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").Value = "'" + '6.57'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").Value = "'" + '11.51'
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("D33").Value = "'" + '580.40'
$ws.Range("E33").Value = '  -9.91%  '
$ws.Range("E34").Value = '  -2.84%  '
$ws.Range("D35").Value = "'" + '59.60'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = "'" + '0.153'
$ws.Range("E36").Value = '  +3.95%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = "'" + '3.53'
$ws.Range("E38").Value = '  +4.01%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = "'" + '36.13'
$ws.Range("E39").Value = '  -6.41%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0767'
$ws.Range("E40").Value = '  -5.71%  '
$ws.Range("E41").Value = '  -4.17%  '
$ws.Range("D42").Value = '3.184.29'
$ws.Range("E42").Value = '  +4.58%  '
$ws.Range("D43").Value = "'" + '1.00'
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("E46").Value = '  -6.00%  '
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  -4.81%  '
$ws.Range("D50").Value = "'" + '8.43'
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("D51").Value = "'" + '137.06'
$ws.Range("E51").Value = '  -3.54%  '
